$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testresults")
$ws.Activate()

# Insert a new row at position 26 (shifts existing rows 26-48 down to 27-49),
# inheriting the number formats (date style in col A, 0.0 style in col H)
# from the row that used to occupy row 26.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new test-result entry.
$ws.Cells.Item(26, 1).Value = 43302
$ws.Cells.Item(26, 2).Value = "MS 20.7.18"
$ws.Cells.Item(26, 3).Value = 11336459
$ws.Cells.Item(26, 4).Value = "ECE-X3"
$ws.Cells.Item(26, 5).Value = "60+0.6"
$ws.Cells.Item(26, 6).Value = "8moves_v3 rnd."
$ws.Cells.Item(26, 7).Value = 500
$ws.Cells.Item(26, 8).Value = 42.3
$ws.Cells.Item(26, 9).Value = "163-97.240"

# Re-establish the AutoFilter over the grown range A1:J49 (row insert keeps
# AutoFilterMode on, but leaves the filter range pinned to the old extent).
$ws.AutoFilterMode = $false
$ws.Range("A1:J49").AutoFilter()

# The hidden _FilterDatabase defined name also needs to track the new extent.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Testresults!`$A`$1:`$J`$49"

# Match the author's final selection/scroll position.
$ws.Range("C7").Select()
$ws.Range("J26").Select()
